$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ==================================================================
# Convertir las formulas originales (Tabla 1 y Tabla 2) en formulas
# compartidas, para que las formulas nuevas compartan el mismo grupo
# si="0", si="1", si="2", si="3" en orden.
# ==================================================================
$ws.Range("G2:G11").Formula = "=AVERAGE(B2:F2)"
$ws.Range("G14:G23").Formula = "=AVERAGE(B14:F14)"

# ==================================================================
# Tabla 3 (filas 25-36): copia de la Tabla 1 (filas 1-11)
# ==================================================================
$ws.Range("A25").Value = "Codigo compañero"

# Copiar el formato (bordes) de la Tabla 1 hacia las nuevas filas B27:G36
$ws.Range("B2:G11").Copy()
$ws.Range("B27").PasteSpecial(-4122)

# Encabezados fila 26 (idéntico a la fila 1)
$ws.Range("A1:G1").Copy()
$ws.Range("A26").PasteSpecial(-4163)
$ws.Range("I1").Copy()
$ws.Range("I26").PasteSpecial(-4163)

$ws.Range("A27").Value = 324
$ws.Range("B27").Value = 106900
$ws.Range("C27").Value = 107500
$ws.Range("D27").Value = 107500
$ws.Range("E27").Value = 106400
$ws.Range("F27").Value = 171400
$ws.Range("I2").Copy()
$ws.Range("I27").PasteSpecial(-4163)
$ws.Range("A28").Value = 43423
$ws.Range("B28").Value = 495300
$ws.Range("C28").Value = 495000
$ws.Range("D28").Value = 432400
$ws.Range("E28").Value = 498100
$ws.Range("F28").Value = 503200
$ws.Range("I3").Copy()
$ws.Range("I28").PasteSpecial(-4163)
$ws.Range("A29").Value = 235356
$ws.Range("B29").Value = 2211500
$ws.Range("C29").Value = 2209000
$ws.Range("D29").Value = 2206900
$ws.Range("E29").Value = 2228100
$ws.Range("F29").Value = 2209700
$ws.Range("I4").Copy()
$ws.Range("I29").PasteSpecial(-4163)
$ws.Range("A30").Value = 234324
$ws.Range("B30").Value = 2206500
$ws.Range("C30").Value = 2203500
$ws.Range("D30").Value = 1576600
$ws.Range("E30").Value = 2231300
$ws.Range("F30").Value = 2205000
$ws.Range("A31").Value = 5654774
$ws.Range("B31").Value = 10027600
$ws.Range("C31").Value = 8015100
$ws.Range("D31").Value = 9489800
$ws.Range("E31").Value = 12047400
$ws.Range("F31").Value = 9545200
$ws.Range("A32").Value = 5235124
$ws.Range("B32").Value = 10734400
$ws.Range("C32").Value = 11586500
$ws.Range("D32").Value = 10644100
$ws.Range("E32").Value = 12594400
$ws.Range("F32").Value = 14598100
$ws.Range("A33").Value = 21213213
$ws.Range("B33").Value = 26071500
$ws.Range("C33").Value = 27385000
$ws.Range("D33").Value = 26859500
$ws.Range("E33").Value = 25500600
$ws.Range("F33").Value = 26957400
$ws.Range("A34").Value = 56863435
$ws.Range("B34").Value = 62803800
$ws.Range("C34").Value = 62790100
$ws.Range("D34").Value = 64514800
$ws.Range("E34").Value = 61111600
$ws.Range("F34").Value = 71477400
$ws.Range("A35").Value = 436346436
$ws.Range("B35").Value = 444088700
$ws.Range("C35").Value = 439661300
$ws.Range("D35").Value = 436972400
$ws.Range("E35").Value = 446940100
$ws.Range("F35").Value = 438637300
$ws.Range("A36").Value = 565235323
$ws.Range("B36").Value = 543234200
$ws.Range("C36").Value = 543758000
$ws.Range("D36").Value = 549920600
$ws.Range("E36").Value = 557457400
$ws.Range("F36").Value = 542130600

$ws.Range("G27:G36").Formula = "=AVERAGE(B27:F27)"

# ==================================================================
# Tabla 4 (filas 38-49): copia de la Tabla 2 (filas 13-23)
# ==================================================================
$ws.Range("A38").Value = "Codigo compañero"

# Copiar el formato (bordes) de la Tabla 2 hacia las nuevas filas B40:G49
$ws.Range("B14:G23").Copy()
$ws.Range("B40").PasteSpecial(-4122)

# Encabezados fila 39 (idéntico a la fila 1, misma columna I que Tabla 1)
$ws.Range("A1:G1").Copy()
$ws.Range("A39").PasteSpecial(-4163)
$ws.Range("I1").Copy()
$ws.Range("I39").PasteSpecial(-4163)

$ws.Range("A40").Value = 324
$ws.Range("B40").Value = 113800
$ws.Range("C40").Value = 107800
$ws.Range("D40").Value = 107200
$ws.Range("E40").Value = 107200
$ws.Range("F40").Value = 107900
$ws.Range("I2").Copy()
$ws.Range("I40").PasteSpecial(-4163)
$ws.Range("A41").Value = 43423
$ws.Range("B41").Value = 350000
$ws.Range("C41").Value = 496300
$ws.Range("D41").Value = 490700
$ws.Range("E41").Value = 349800
$ws.Range("F41").Value = 495800
$ws.Range("I3").Copy()
$ws.Range("I41").PasteSpecial(-4163)
$ws.Range("A42").Value = 235356
$ws.Range("B42").Value = 2213100
$ws.Range("C42").Value = 2213600
$ws.Range("D42").Value = 2209600
$ws.Range("E42").Value = 2070300
$ws.Range("F42").Value = 2217500
$ws.Range("I4").Copy()
$ws.Range("I42").PasteSpecial(-4163)
$ws.Range("A43").Value = 234324
$ws.Range("B43").Value = 2206800
$ws.Range("C43").Value = 1529200
$ws.Range("D43").Value = 2205000
$ws.Range("E43").Value = 2204600
$ws.Range("F43").Value = 2188100
$ws.Range("A44").Value = 5654774
$ws.Range("B44").Value = 9557400
$ws.Range("C44").Value = 9632500
$ws.Range("D44").Value = 9471600
$ws.Range("E44").Value = 8026900
$ws.Range("F44").Value = 9587900
$ws.Range("A45").Value = 5235124
$ws.Range("B45").Value = 9218700
$ws.Range("C45").Value = 10774500
$ws.Range("D45").Value = 9158700
$ws.Range("E45").Value = 10814800
$ws.Range("F45").Value = 11239400
$ws.Range("A46").Value = 21213213
$ws.Range("B46").Value = 33986600
$ws.Range("C46").Value = 31427400
$ws.Range("D46").Value = 27271300
$ws.Range("E46").Value = 25931500
$ws.Range("F46").Value = 26964400
$ws.Range("A47").Value = 56863435
$ws.Range("B47").Value = 79262600
$ws.Range("C47").Value = 69104400
$ws.Range("D47").Value = 69605200
$ws.Range("E47").Value = 61927400
$ws.Range("F47").Value = 65044800
$ws.Range("A48").Value = 436346436
$ws.Range("B48").Value = 442941900
$ws.Range("C48").Value = 445144600
$ws.Range("D48").Value = 439727300
$ws.Range("E48").Value = 442382500
$ws.Range("F48").Value = 439424200
$ws.Range("A49").Value = 565235323
$ws.Range("B49").Value = 553828300
$ws.Range("C49").Value = 544060500
$ws.Range("D49").Value = 551903100
$ws.Range("E49").Value = 542087000
$ws.Range("F49").Value = 549775400

$ws.Range("G40:G49").Formula = "=AVERAGE(B40:F40)"

# ==================================================================
# Vista de la hoja: selección final y celda superior izquierda
# ==================================================================
$ws.Range("F51").Select()
